$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append four new rows (18-21) of test-result data to the sheet,
# mirroring the existing pattern of Scenario/Result/Browser columns.

$ws.Range("A18").Value = "Login with valid username and password"
$ws.Range("B18").Value = "FAILED"
$ws.Range("C18").Value = "chrome"

$ws.Range("A19").Value = "Login with valid username and password"
$ws.Range("B19").Value = "FAILED"
$ws.Range("C19").Value = "chrome"

$ws.Range("A20").Value = "Login with valid username and password"
$ws.Range("B20").Value = "PASSED"
$ws.Range("C20").Value = "chrome"

$ws.Range("A21").Value = "Create a country"
$ws.Range("B21").Value = "PASSED"
$ws.Range("C21").Value = "chrome"
